$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.36437463760376
$ws.Range("B1").Value = 3.397713661193848
$ws.Range("C1").Value = 3.030774831771851
$ws.Range("D1").Value = 2.515383243560791
$ws.Range("E1").Value = 1.642297387123108
